# This script reproduces the pipeline re-run update described in the commit
# "Added and updated mounted pipeline": the underlying per-segment signal
# histograms (Step1_Data), their cumulative-sum normalizations (Step2_Sj) and the
# derived threshold crossing statistics (Step3_DataPts_0.5/0.7/0.8/0.9) were
# recomputed upstream and the refreshed values are written back into the cells
# below (values only; no cells are added, removed, or restyled).

$wb = $excel.ActiveWorkbook

# --- Step1_Data ---
$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0.2127069668866213
$ws.Range("F2").Value = 0.09863084541182979
$ws.Range("G2").Value = 0.02240959447844173
$ws.Range("I2").Value = 0.08870159263506186
$ws.Range("J2").Value = 0.002834380685868628
$ws.Range("K2").Value = 0.01047825558168749
$ws.Range("L2").Value = 0.004995991788587114
$ws.Range("M2").Value = 0.1617046486310822
$ws.Range("N2").Value = 0.06104055130028296
$ws.Range("O2").Value = 0.004614461453189832
$ws.Range("P2").Value = 0.00201803711090926
$ws.Range("T2").Value = 0.06757921243725445
$ws.Range("U2").Value = 0.04099505719257693
$ws.Range("V2").Value = 0.02274923940574862
$ws.Range("W2").Value = 0.002627297750365359
$ws.Range("X2").Value = 0.009587423839133254
$ws.Range("Y2").Value = 0.001393273093452157
$ws.Range("Z2").Value = 0.005764571114115846
$ws.Range("AA2").Value = 0.01275205596089575
$ws.Range("AB2").Value = 0.01483437447993228
$ws.Range("AC2").Value = 0.0006946681225420735
$ws.Range("AD2").Value = 0.05962983970657044
$ws.Range("AE2").Value = 0.04106151381231295
$ws.Range("AF2").Value = 0.01660536206909148
$ws.Range("AG2").Value = 0.02585628759214127
$ws.Range("AH2").Value = 0.006207537190651822
$ws.Range("AJ2").Value = 0.001526960269653014
$ws.Range("F3").Value = 0.003589283571538827
$ws.Range("J3").Value = 0.004894497176005506
$ws.Range("M3").Value = 0.2096146787542453
$ws.Range("N3").Value = 0.01234425148966898
$ws.Range("O3").Value = 0.04440259957954138
$ws.Range("T3").Value = 0.06235320421401432
$ws.Range("V3").Value = 0.02684863933364152
$ws.Range("AF3").Value = 0.03509823116582758
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.1235540490229708
$ws.Range("F5").Value = 0.242091828412845
$ws.Range("G5").Value = 0.01063765458933459
$ws.Range("I5").Value = 0.0194069277310621
$ws.Range("K5").Value = 0.02282299073868249
$ws.Range("L5").Value = 0.06004126313402935
$ws.Range("M5").Value = 0.00281715681312533
$ws.Range("N5").Value = 0.1154592514719948
$ws.Range("P5").Value = 0.003958071575513147
$ws.Range("S5").Value = 0.006554720145677944
$ws.Range("T5").Value = 0.07703755903486988
$ws.Range("U5").Value = 0.02546605472640359
$ws.Range("V5").Value = 0.02980550791041364
$ws.Range("X5").Value = 0.01378415922730212
$ws.Range("Z5").Value = 0.005461746682350133
$ws.Range("AA5").Value = 0.02271565561168362
$ws.Range("AC5").Value = 0.003162714641925902
$ws.Range("AD5").Value = 0.1149599576381318
$ws.Range("AE5").Value = 0.0537466767057477
$ws.Range("AG5").Value = 0.01308213168038198
$ws.Range("AI5").Value = 0.03343392250555424
$ws.Range("E6").Value = 0.1867706865320108
$ws.Range("F6").Value = 0.09549714414984097
$ws.Range("G6").Value = 0.003451895100600371
$ws.Range("I6").Value = 0.04450690356453555
$ws.Range("K6").Value = 0.001006741114362826
$ws.Range("L6").Value = 0.05620356835582065
$ws.Range("M6").Value = 0.2199572111014031
$ws.Range("N6").Value = 0.002957731339256682
$ws.Range("O6").Value = 0.002442741628188215
$ws.Range("T6").Value = 0.1379820371023241
$ws.Range("U6").Value = 0.03093582241426343
$ws.Range("V6").Value = 0.008995219297606495
$ws.Range("X6").Value = 0.03661983945777809
$ws.Range("Z6").Value = 0.005065952712430966
$ws.Range("AA6").Value = 0.01047751565178345
$ws.Range("AB6").Value = 0.008456001261579549
$ws.Range("AC6").Value = 0.001620549746577602
$ws.Range("AD6").Value = 0.07064339767375213
$ws.Range("AE6").Value = 0.03975243500873663
$ws.Range("AF6").Value = 0.0192344437176266
$ws.Range("AG6").Value = 0.007423015390198904
$ws.Range("AH6").Value = 0.006630684982795576
$ws.Range("AJ6").Value = 0.003368462696527795
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0.1874767773268845
$ws.Range("F11").Value = 0.1262590159419168
$ws.Range("G11").Value = 0.07615149604906249
$ws.Range("H11").Value = 0.0112157930349232
$ws.Range("I11").Value = 0.008263472450238804
$ws.Range("J11").Value = 0.003301417528353657
$ws.Range("K11").Value = 0.007682651997925804
$ws.Range("L11").Value = 0.1228605233532873
$ws.Range("M11").Value = 0.01175396962012543
$ws.Range("N11").Value = 0.06094950138474729
$ws.Range("O11").Value = 0.0766426800260539
$ws.Range("P11").Value = 0.01008211862815695
$ws.Range("T11").Value = 0.05171517489825327
$ws.Range("U11").Value = 0.0572242084104311
$ws.Range("V11").Value = 0.01359695373223934
$ws.Range("X11").Value = 0.00375802268334008
$ws.Range("Y11").Value = 0.005336754379508112
$ws.Range("Z11").Value = 0.000919475117589948
$ws.Range("AA11").Value = 0.01096747962367951
$ws.Range("AC11").Value = 0.006259812806766983
$ws.Range("AD11").Value = 0.06282068488938651
$ws.Range("AE11").Value = 0.06179785921607003
$ws.Range("AF11").Value = 0.004635531801567031
$ws.Range("AI11").Value = 0.005356893143387767
$ws.Range("AJ11").Value = 0.01297173195610424

# --- Step2_Sj ---
$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0.2127069668866213
$ws.Range("F2").Value = 0.3113378122984511
$ws.Range("G2").Value = 0.3337474067768929
$ws.Range("H2").Value = 0.3337474067768929
$ws.Range("I2").Value = 0.4224489994119547
$ws.Range("J2").Value = 0.4252833800978233
$ws.Range("K2").Value = 0.4357616356795108
$ws.Range("L2").Value = 0.4407576274680979
$ws.Range("M2").Value = 0.6024622760991801
$ws.Range("N2").Value = 0.6635028273994631
$ws.Range("O2").Value = 0.6681172888526529
$ws.Range("P2").Value = 0.6701353259635622
$ws.Range("Q2").Value = 0.6701353259635622
$ws.Range("R2").Value = 0.6701353259635622
$ws.Range("S2").Value = 0.6701353259635622
$ws.Range("T2").Value = 0.7377145384008166
$ws.Range("U2").Value = 0.7787095955933935
$ws.Range("V2").Value = 0.8014588349991421
$ws.Range("W2").Value = 0.8040861327495075
$ws.Range("X2").Value = 0.8136735565886407
$ws.Range("Y2").Value = 0.8150668296820929
$ws.Range("Z2").Value = 0.8208314007962088
$ws.Range("AA2").Value = 0.8335834567571045
$ws.Range("AB2").Value = 0.8484178312370368
$ws.Range("AC2").Value = 0.8491124993595789
$ws.Range("AD2").Value = 0.9087423390661493
$ws.Range("AE2").Value = 0.9498038528784623
$ws.Range("AF2").Value = 0.9664092149475537
$ws.Range("AG2").Value = 0.992265502539695
$ws.Range("AH2").Value = 0.9984730397303468
$ws.Range("AI2").Value = 0.9984730397303468
$ws.Range("AJ2").Value = 0.9999999999999999
$ws.Range("H3").Value = 0.2499601219552051
$ws.Range("AE3").Value = 0.9012657964565257
$ws.Range("AF3").Value = 0.9363640276223533
$ws.Range("AG3").Value = 0.9851456292440508
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.1235540490229708
$ws.Range("F5").Value = 0.3656458774358158
$ws.Range("G5").Value = 0.3762835320251504
$ws.Range("H5").Value = 0.3762835320251504
$ws.Range("I5").Value = 0.3956904597562125
$ws.Range("J5").Value = 0.3956904597562125
$ws.Range("K5").Value = 0.418513450494895
$ws.Range("L5").Value = 0.4785547136289243
$ws.Range("M5").Value = 0.4813718704420497
$ws.Range("N5").Value = 0.5968311219140445
$ws.Range("O5").Value = 0.5968311219140445
$ws.Range("P5").Value = 0.6007891934895577
$ws.Range("Q5").Value = 0.6007891934895577
$ws.Range("R5").Value = 0.6007891934895577
$ws.Range("S5").Value = 0.6073439136352357
$ws.Range("T5").Value = 0.6843814726701055
$ws.Range("U5").Value = 0.7098475273965091
$ws.Range("V5").Value = 0.7396530353069227
$ws.Range("W5").Value = 0.7396530353069227
$ws.Range("X5").Value = 0.7534371945342249
$ws.Range("Y5").Value = 0.7534371945342249
$ws.Range("Z5").Value = 0.758898941216575
$ws.Range("AA5").Value = 0.7816145968282586
$ws.Range("AB5").Value = 0.7816145968282586
$ws.Range("AC5").Value = 0.7847773114701845
$ws.Range("AD5").Value = 0.8997372691083163
$ws.Range("AE5").Value = 0.953483945814064
$ws.Range("AF5").Value = 0.953483945814064
$ws.Range("AG5").Value = 0.966566077494446
$ws.Range("AH5").Value = 0.966566077494446
$ws.Range("E6").Value = 0.1867706865320108
$ws.Range("F6").Value = 0.2822678306818517
$ws.Range("G6").Value = 0.2857197257824521
$ws.Range("H6").Value = 0.2857197257824521
$ws.Range("I6").Value = 0.3302266293469877
$ws.Range("J6").Value = 0.3302266293469877
$ws.Range("K6").Value = 0.3312333704613505
$ws.Range("L6").Value = 0.3874369388171711
$ws.Range("M6").Value = 0.6073941499185742
$ws.Range("N6").Value = 0.6103518812578309
$ws.Range("O6").Value = 0.6127946228860192
$ws.Range("P6").Value = 0.6127946228860192
$ws.Range("Q6").Value = 0.6127946228860192
$ws.Range("R6").Value = 0.6127946228860192
$ws.Range("S6").Value = 0.6127946228860192
$ws.Range("T6").Value = 0.7507766599883432
$ws.Range("U6").Value = 0.7817124824026066
$ws.Range("V6").Value = 0.7907077017002131
$ws.Range("W6").Value = 0.7907077017002131
$ws.Range("X6").Value = 0.8273275411579912
$ws.Range("Y6").Value = 0.8273275411579912
$ws.Range("Z6").Value = 0.8323934938704222
$ws.Range("AA6").Value = 0.8428710095222056
$ws.Range("AB6").Value = 0.8513270107837851
$ws.Range("AC6").Value = 0.8529475605303627
$ws.Range("AD6").Value = 0.9235909582041149
$ws.Range("AE6").Value = 0.9633433932128515
$ws.Range("AF6").Value = 0.9825778369304781
$ws.Range("AG6").Value = 0.9900008523206769
$ws.Range("AH6").Value = 0.9966315373034725
$ws.Range("AI6").Value = 0.9966315373034725
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0.1874767773268845
$ws.Range("F11").Value = 0.3137357932688012
$ws.Range("G11").Value = 0.3898872893178638
$ws.Range("H11").Value = 0.401103082352787
$ws.Range("I11").Value = 0.4093665548030258
$ws.Range("J11").Value = 0.4126679723313794
$ws.Range("K11").Value = 0.4203506243293052
$ws.Range("L11").Value = 0.5432111476825925
$ws.Range("M11").Value = 0.5549651173027179
$ws.Range("N11").Value = 0.6159146186874652
$ws.Range("O11").Value = 0.6925572987135191
$ws.Range("P11").Value = 0.7026394173416761
$ws.Range("Q11").Value = 0.7026394173416761
$ws.Range("R11").Value = 0.7026394173416761
$ws.Range("S11").Value = 0.7026394173416761
$ws.Range("T11").Value = 0.7543545922399294
$ws.Range("U11").Value = 0.8115788006503605
$ws.Range("V11").Value = 0.8251757543825998
$ws.Range("W11").Value = 0.8251757543825998
$ws.Range("X11").Value = 0.82893377706594
$ws.Range("Y11").Value = 0.834270531445448
$ws.Range("Z11").Value = 0.835190006563038
$ws.Range("AA11").Value = 0.8461574861867175
$ws.Range("AB11").Value = 0.8461574861867175
$ws.Range("AC11").Value = 0.8524172989934845
$ws.Range("AD11").Value = 0.915237983882871
$ws.Range("AE11").Value = 0.9770358430989411
$ws.Range("AF11").Value = 0.9816713749005082
$ws.Range("AG11").Value = 0.9816713749005082
$ws.Range("AH11").Value = 0.9816713749005082
$ws.Range("AI11").Value = 0.9870282680438959

# --- Step3_DataPts_0.5 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F2").Value = 0.6024622760991801
$ws.Range("F5").Value = 0.5968311219140445
$ws.Range("F6").Value = 0.6073941499185742
$ws.Range("F11").Value = 0.5432111476825925

# --- Step3_DataPts_0.7 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("F2").Value = 0.7377145384008166
$ws.Range("F5").Value = 0.7098475273965091
$ws.Range("F6").Value = 0.7507766599883432
$ws.Range("D11").Value = 15
$ws.Range("F11").Value = 0.7026394173416761
$ws.Range("G11").Value = 13

# --- Step3_DataPts_0.8 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("F2").Value = 0.8014588349991421
$ws.Range("F5").Value = 0.8997372691083163
$ws.Range("F6").Value = 0.8273275411579912
$ws.Range("F11").Value = 0.8115788006503605

# --- Step3_DataPts_0.9 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("F2").Value = 0.9087423390661493
$ws.Range("F3").Value = 0.9012657964565257
$ws.Range("D5").Value = 30
$ws.Range("F5").Value = 0.953483945814064
$ws.Range("G5").Value = 28
$ws.Range("F6").Value = 0.9235909582041149
$ws.Range("F11").Value = 0.915237983882871
